$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8 block: "By Register" -> "Direct By Register" ---
$ws.Range("A8").Value = "Direct By Register"

# --- Row 14/15 block: "Direct By Register" -> "Indirect By Register" ---
$ws.Range("A14").Value = "Indirect By Register"

# Row 14: was D14=[RR3] Hi, E14=[RR3] Lo, (F14:G14 merged)=NOT USED
# Becomes: (D14:G14 merged)=NOT USED, new H14 = "Address is implict in H and L registers"
$ws.Range("F14:G14").UnMerge()
$ws.Range("D14").Value = "NOT USED"
$ws.Range("E14").ClearContents()
$ws.Range("F14").ClearContents()
$ws.Range("G14").ClearContents()

# Copy cell formatting (borders/font) from the equivalent row-8 header cells so the
# merged D14:G14 block keeps the same "left / middle / middle / right" bordered look.
$ws.Range("E8").Copy()
$ws.Range("D14").PasteSpecial(-4122)
$ws.Range("F8").Copy()
$ws.Range("E14").PasteSpecial(-4122)
$ws.Range("F8").Copy()
$ws.Range("F14").PasteSpecial(-4122)
$ws.Range("G8").Copy()
$ws.Range("G14").PasteSpecial(-4122)

$ws.Range("D14:G14").Merge()

$ws.Range("H14").Value = "Address is implict in H and L registers"

# Row 15: was D15=3 bits, E15=3 bits, (F15:G15 merged)=9 bits
# Becomes: (D15:G15 merged)=15 bits
$ws.Range("F15:G15").UnMerge()
$ws.Range("D15").Value = "15 bits"
$ws.Range("E15").ClearContents()
$ws.Range("F15").ClearContents()
$ws.Range("G15").ClearContents()

$ws.Range("E9").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("F9").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("F9").Copy()
$ws.Range("F15").PasteSpecial(-4122)
$ws.Range("G9").Copy()
$ws.Range("G15").PasteSpecial(-4122)

$ws.Range("D15:G15").Merge()

# --- Selection moved from I6 to H15 ---
$ws.Range("H15").Select()
